$wb = $excel.ActiveWorkbook

$wsDoc = $wb.Worksheets.Item("DOCUMENT")
$wsCov = $wb.Worksheets.Item("COVERAGE")

# Replace the fully-qualified ConcatAggregator class reference with the
# shorter "Concat" alias throughout the JETT template cells.
$wsDoc.Range("B6").Value = "<jt:forEach items=`"`${source.requirements}`" var=`"req`"><jt:if test=`"`${req.referredBySource.isEmpty() && !source.coversBy.isEmpty()}`" elseAction=`"shiftUp`">`${jagg:eval(req.referredBySource, 'Concat(name, `"'+line+'`")')}</jt:if>"
$wsDoc.Range("B7").Value = "<jt:if test=`"`${!req.referredBySource.isEmpty() && req.referredBySource.size()<source.coversBy.size()}`" elseAction=`"shiftUp`">`${jagg:eval(req.referredBySource, 'Concat(name, `"'+line+'`")')}</jt:if>"
$wsDoc.Range("B8").Value = "<jt:if test=`"`${req.referredBySource.size()>=source.coversBy.size()}`" elseAction=`"shiftUp`">`${jagg:eval(req.referredBySource, 'Concat(name, `"'+line+'`")')}</jt:if></jt:forEach>"

$wsCov.Range("B6").Value = "`${jagg:eval(req.getReferencesFor(coverSource), 'Concat(text, `"'+line+'`")')}</jt:if></jt:forEach>"
$wsCov.Range("B10").Value = "`${jagg:eval(req.getReferredByRequirementFor(coverBySource), 'Concat(text, `"'+line+'`")')}</jt:if></jt:forEach>"

# Move the active selection / active sheet: the author was last looking at
# the COVERAGE sheet, with the cursor on B23 there, while the DOCUMENT
# sheet cursor sits on B8.
$wsDoc.Range("B8").Select() | Out-Null
$wsCov.Activate()
$wsCov.Range("B23").Select() | Out-Null
